$wb = $excel.ActiveWorkbook

# --- Fix the mislabeled 2050 column header (E1) on the tables that show it. ---
# A plain Range("E1").Value = "2050" gets auto-coerced to the NUMBER 2050 by
# Excel (it "looks like a number"), which is not what we want - the header
# must stay a literal text label, exactly like the other year headers
# (B1/C1/D1) on the same row. To force a text value while preserving the
# cell's existing style (bold/border/center), write the literal text into a
# scratch cell via a formula (="2050") - which makes it a text-typed
# formula result - then Copy/PasteSpecial "values only" into E1. That
# carries over only the resulting text value, not the scratch cell's
# formula or formatting, so E1 keeps its original style untouched.
function Set-TextHeader {
    param($ws, [string]$cellAddr, [string]$text)

    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

$sheetNames2050 = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $sheetNames2050) {
    $ws = $wb.Worksheets.Item($name)
    Set-TextHeader $ws "E1" "2050"
}

# Sheet 4 uses year-range headers ("2015-2030", "2031-2040"), so the last
# column's header should read "2041-2050" (not a pure-numeric string, so no
# auto-coercion risk - but use the same helper for consistency).
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
Set-TextHeader $ws4 "E1" "2041-2050"

# --- Remove the "Total" row (row 13) from the four detailed tables. ---
$sheetsWithTotalRow13 = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($name in $sheetsWithTotalRow13) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(13).Delete()
}

# --- Remove the "Total" row (row 4) from the cost table. ---
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Rows.Item(4).Delete()
